$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 50 (old rows 50-77 shift down to 52-79)
$ws.Rows.Item(50).Insert()
$ws.Rows.Item(50).Insert()

$d = (Get-Date -Year 2022 -Month 11 -Day 10 -Hour 0 -Minute 0 -Second 0).Date

# New row 50
$ws.Range("A50").Value = 10
$ws.Range("B50").Value = "Vega Modelo de Temuco"
$ws.Range("C50").Value = "La Araucanía"
$ws.Range("D50").Value = $d
$ws.Range("E50").Value = 9
$ws.Range("F50").Value = 300000000
$ws.Range("G50").Value = "Espárragos"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 500
$ws.Range("K50").Value = 1500
$ws.Range("L50").Value = 1500
$ws.Range("M50").Value = 1500
$ws.Range("N50").Value = "$/kilo"
$ws.Range("O50").Value = "Región de La Araucanía"
$ws.Range("P50").Value = 1500
$ws.Range("Q50").Value = 1
$ws.Range("R50").Value = "Hortaliza"

# New row 51
$ws.Range("A51").Value = 10
$ws.Range("B51").Value = "Vega Modelo de Temuco"
$ws.Range("C51").Value = "La Araucanía"
$ws.Range("D51").Value = $d
$ws.Range("E51").Value = 9
$ws.Range("F51").Value = 300000000
$ws.Range("G51").Value = "Espárragos"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 300
$ws.Range("K51").Value = 1500
$ws.Range("L51").Value = 1500
$ws.Range("M51").Value = 1500
$ws.Range("N51").Value = "$/kilo"
$ws.Range("O51").Value = "Región del Maule"
$ws.Range("P51").Value = 1500
$ws.Range("Q51").Value = 1
$ws.Range("R51").Value = "Hortaliza"

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
